$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text values look like plain numbers must be forced to
# Text format first, otherwise Excel auto-converts them to doubles (losing
# formatting such as trailing zeros, e.g. "1.600" -> 1.6) when .Value is set.

# Row 2
$ws.Range("D2").Value = '26.912.27'
$ws.Range("E2").Value = '  +1.18%  '

# Row 3
$ws.Range("D3").Value = '1.846.02'
$ws.Range("E3").Value = '  +1.44%  '

# Row 4
$ws.Range("E4").Value = '  +0.09%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '309.25'
$ws.Range("E5").Value = '  +1.19%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.006'
$ws.Range("E6").Value = '  +0.00%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4806'
$ws.Range("E7").Value = '  +2.99%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3678'
$ws.Range("E8").Value = '  +2.21%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07195'
$ws.Range("E9").Value = '  +1.02%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.9227'
$ws.Range("E10").Value = '  +2.28%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '19.61'
$ws.Range("E11").Value = '  +0.92%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07609'
$ws.Range("E12").Value = '  -2.21%  '

# Row 13
$ws.Range("D13").Value = '1.863.10'
$ws.Range("E13").Value = '  +1.30%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.310'
$ws.Range("E14").Value = '  +1.22%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.397'
$ws.Range("E15").Value = '  +1.00%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '88.35'
$ws.Range("E16").Value = '  +1.01%  '

# Row 17
$ws.Range("E17").Value = '  +0.12%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008651'
$ws.Range("E18").Value = '  +1.10%  '

# Row 19
$ws.Range("E19").Value = '  +0.02%  '

# Row 20
$ws.Range("D20").Value = '26.952.94'
$ws.Range("E20").Value = '  +1.16%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.54'
$ws.Range("E21").Value = '  +2.41%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.035'
$ws.Range("E22").Value = '  +0.52%  '

# Row 23
$ws.Range("E23").Value = '  +0.83%  '

# Row 24
$ws.Range("E24").Value = '  -0.79%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '152.30'
$ws.Range("E25").Value = '  +0.22%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '18.15'
$ws.Range("E26").Value = '  +1.51%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.003'
$ws.Range("E27").Value = '  +1.03%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '114.24'
$ws.Range("E28").Value = '  +0.63%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '4.907'
$ws.Range("E29").Value = '  +1.89%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.08876'
$ws.Range("E30").Value = '  +1.15%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.274'
$ws.Range("E31").Value = '  +4.40%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.7461'
$ws.Range("E32").Value = '  +1.70%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.164'
$ws.Range("E33").Value = '  +3.73%  '

# Row 34
$ws.Range("B34").Value = 'Filecoin'
$ws.Range("C34").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.492'
$ws.Range("E34").Value = '  +1.21%  '

# Row 35
$ws.Range("B35").Value = 'RenderToken'
$ws.Range("C35").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.758'
$ws.Range("E35").Value = '  -0.93%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.089'
$ws.Range("E36").Value = '  +1.20%  '

# Row 37
$ws.Range("B37").Value = 'VeChain'
$ws.Range("C37").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.01949'
$ws.Range("E37").Value = '  +0.89%  '

# Row 38
$ws.Range("B38").Value = 'Hedera'
$ws.Range("C38").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.05258'
$ws.Range("E38").Value = '  +2.74%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.966'
$ws.Range("E39").Value = '  +1.92%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.5210'
$ws.Range("E40").Value = '  +2.86%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '6.896'
$ws.Range("E41").Value = '  +1.41%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1512'
$ws.Range("E42").Value = '  +1.08%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '8.198'
$ws.Range("E43").Value = '  +2.54%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '10.53'
$ws.Range("E44").Value = '  +5.32%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.4695'
$ws.Range("E45").Value = '  +0.30%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.006'
$ws.Range("E46").Value = '  +0.00%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '102.09'
$ws.Range("E47").Value = '  +3.45%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.600'
$ws.Range("E48").Value = '  +2.25%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '65.48'
$ws.Range("E49").Value = '  +2.57%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.06036'
$ws.Range("E50").Value = '  +0.47%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.8852'
$ws.Range("E51").Value = '  +4.23%  '
